$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H), matching the formatting of the existing
# header cells (bold, centered, bordered) by copying G1's format onto it.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data value under the "Save" header.
$ws.Range("H2").Value = 1
